# "update hotel reviews data"
#
# The hotel_info sheet had three review-statistic columns
# (English_Reviews_num, Local_Rank, Total_Reviews_num) that were still
# blank for this hotel's row. Fill them in with the scraped values.
#
# NOTE: in the source data these numeric-looking values are stored as
# text (shared-string) cells, not numbers -- format the range as Text
# first so Excel doesn't auto-coerce "5"/"421" into numeric cells, then
# drop the temporary formatting again so no stray number format is left
# behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$targetRange = $ws.Range("G2:I2")
$targetRange.NumberFormat = "@"

$ws.Range("G2").Value = "5"    # English_Reviews_num
$ws.Range("H2").Value = "421"  # Local_Rank
$ws.Range("I2").Value = "5"    # Total_Reviews_num

$targetRange.ClearFormats()
